# Apply the documentation sheet updates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "DATABASE @" -> "DjATABASE @" ---
$ws.Range("B8").Value = "DjATABASE @"

# --- Row 9: djatabase -> root / r00tc4nrun ---
$ws.Range("B9").Value = "root"
$ws.Range("C9").Value = "r00tc4nrun"

# --- Row 10: root/r00tc4nrun -> djonas/dj0n@s (becomes the hyperlinked row) ---
$ws.Range("B10").Value = "djonas"
$ws.Range("C10").Value = "dj0n@s"

# --- Row 11: djonas/dj0n@s -> php myadmin (loses the hyperlink) ---
$ws.Range("C11").Hyperlinks.Delete() | Out-Null
$ws.Range("B11").Value = "php myadmin"
$ws.Range("C11").Clear()

# Re-create the hyperlink on its new location (C10) with the Hyperlink style.
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:dj0n@s") | Out-Null
$ws.Range("C10").Style = "Hyperlink"

# --- Row 12: php myadmin -> root / r00tc4nrun ---
$ws.Range("B12").Value = "root"
$ws.Range("C12").Value = "r00tc4nrun"

# --- Rows 13 & 14 removed entirely ---
$ws.Range("B13:C13").Clear()
$ws.Range("B14:C14").Clear()

# --- Update the saved selection to B13 ---
$ws.Range("B13").Select() | Out-Null
